$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------------
# 1. Move the task that was sitting in "To do" (A25) into the "Done" column
#    (C25), matching the existing formatting used for that area of the table.
# ---------------------------------------------------------------------------
$movedText = $ws.Range("A25").Value2
$ws.Range("A25").ClearContents()

# ---------------------------------------------------------------------------
# 2. Grow the table so the new tasks become real table rows (this keeps the
#    table ref / autofilter / dimension consistent with the new data).
# ---------------------------------------------------------------------------
$tbl.Resize($ws.Range("A1:C30"))

# ---------------------------------------------------------------------------
# 3. Re-apply the row formatting that the table used for these rows before
#    the resize (PasteSpecial formats only, values are set afterwards).
# ---------------------------------------------------------------------------
$ws.Range("A20:B20").Copy()
$ws.Range("A25:B25").PasteSpecial(-4122)
$ws.Range("A27:B27").PasteSpecial(-4122)
$ws.Range("A28:B28").PasteSpecial(-4122)
$ws.Range("A29:B29").PasteSpecial(-4122)
$ws.Range("A30:B30").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C30").PasteSpecial(-4122)

$ws.Range("A19:B19").Copy()
$ws.Range("A26:B26").PasteSpecial(-4122)
$ws.Range("B19").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Write the cell values: the moved "Done" task plus the 5 new tasks that
#    were added to the "Done" column.
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = $movedText
$ws.Range("C26").Value = "Configurar React y organizar la estructura del proyecto"
$ws.Range("C27").Value = "Crear una página de inicio básica"
$ws.Range("C28").Value = "Configurar CORS en Django para permitir conexiones desde React"
$ws.Range("C29").Value = "Crear servicio para conectar con la API"
$ws.Range("C30").Value = "Probar conexión entre React y Django"

# ---------------------------------------------------------------------------
# 5. Widen columns B:C slightly to fit the new (longer) text.
# ---------------------------------------------------------------------------
$ws.Columns("B:C").ColumnWidth = 54.830729166666664

# ---------------------------------------------------------------------------
# 6. Update the view: drop the old scroll position / selection, zoom to 83%
#    and select F23 (matching the author's final on-screen state).
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 83
[void]$ws.Range("F23").Select()
